$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'40.822.94"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -2.25%  '
$ws.Range("D3").Value = "'2.386.15"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -3.69%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = "'314.32"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.44%  '
$ws.Range("D6").Value = "'88.35"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -5.20%  '
$ws.Range("E7").Value = '  -4.12%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("E9").Value = '  -4.75%  '
$ws.Range("D10").Value = "'0.0825"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -5.17%  '
$ws.Range("D11").Value = "'31.24"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -6.04%  '
$ws.Range("E12").Value = '  -1.85%  '
$ws.Range("D13").Value = "'2.753.44"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -3.74%  '
$ws.Range("D14").Value = "'6.58"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -4.82%  '
$ws.Range("D15").Value = "'15.14"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -3.96%  '
$ws.Range("D16").Value = "'2.356.72"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -4.80%  '
$ws.Range("D17").Value = "'0.762"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -3.92%  '
$ws.Range("D18").Value = "'40.738.09"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -2.38%  '
$ws.Range("D19").Value = "'0.0₃0914"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -3.98%  '
$ws.Range("E20").Value = '  -4.69%  '
$ws.Range("D21").Value = "'69.31"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -2.76%  '
$ws.Range("E22").Value = '  -4.49%  '
$ws.Range("D23").Value = "'232.96"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -2.84%  '
$ws.Range("E24").Value = '  -3.65%  '
$ws.Range("E25").Value = '  +0.08%  '
$ws.Range("D26").Value = "'1.83"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -6.13%  '
$ws.Range("D27").Value = "'23.71"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -4.18%  '
$ws.Range("E28").Value = '  -3.70%  '
$ws.Range("D29").Value = "'9.39"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -4.46%  '
$ws.Range("D30").Value = "'34.02"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -5.93%  '
$ws.Range("D31").Value = "'156.54"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.16%  '
$ws.Range("E32").Value = '  +0.09%  '
$ws.Range("D33").Value = "'5.23"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -5.41%  '
$ws.Range("D34").Value = "'0.0733"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -4.56%  '
$ws.Range("E35").Value = '  -6.69%  '
$ws.Range("E36").Value = '  -2.06%  '
$ws.Range("E37").Value = '  -3.95%  '
$ws.Range("D38").Value = "'16.09"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -8.05%  '
$ws.Range("D39").Value = "'0.100"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -3.01%  '
$ws.Range("E40").Value = '  -7.67%  '
$ws.Range("E41").Value = '  -5.83%  '
$ws.Range("E42").Value = '  -7.84%  '
$ws.Range("D43").Value = "'1.961.94"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.58%  '
$ws.Range("D44").Value = "'0.0271"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -5.18%  '
$ws.Range("D45").Value = "'17.65"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -7.10%  '
$ws.Range("D46").Value = "'2.81"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -6.53%  '
$ws.Range("D47").Value = "'9.34"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.36%  '
$ws.Range("D48").Value = "'2.615.87"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -3.73%  '
$ws.Range("D49").Value = "'93.87"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -3.87%  '
$ws.Range("D50").Value = "'72.92"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.79%  '
$ws.Range("D51").Value = "'51.16"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -3.09%  '
